$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($range, $value)
    $escaped = $value -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

Set-CellText $ws.Range('D2') '41.528.03'
Set-CellText $ws.Range('E2') '  +0.62%  '
Set-CellText $ws.Range('D3') '2.484.32'
Set-CellText $ws.Range('D4') '0.998'
Set-CellText $ws.Range('E4') '  -0.12%  '
Set-CellText $ws.Range('D5') '311.25'
Set-CellText $ws.Range('E5') '  +0.17%  '
Set-CellText $ws.Range('E6') '  -1.72%  '
Set-CellText $ws.Range('E7') '  -1.58%  '
Set-CellText $ws.Range('D8') '1.00'
Set-CellText $ws.Range('E8') '  -0.26%  '
Set-CellText $ws.Range('D9') '0.496'
Set-CellText $ws.Range('E9') '  -2.12%  '
Set-CellText $ws.Range('D10') '32.35'
Set-CellText $ws.Range('E10') '  -3.64%  '
Set-CellText $ws.Range('D11') '0.0780'
Set-CellText $ws.Range('E11') '  -0.17%  '
Set-CellText $ws.Range('D12') '0.110'
Set-CellText $ws.Range('E12') '  +1.36%  '
Set-CellText $ws.Range('D13') '2.870.24'
Set-CellText $ws.Range('E13') '  +0.84%  '
Set-CellText $ws.Range('E14') '  -2.08%  '
Set-CellText $ws.Range('D15') '15.34'
Set-CellText $ws.Range('E15') '  +5.48%  '
Set-CellText $ws.Range('D16') '2.483.21'
Set-CellText $ws.Range('E16') '  +0.98%  '
Set-CellText $ws.Range('D17') '0.759'
Set-CellText $ws.Range('E17') '  -3.81%  '
Set-CellText $ws.Range('D18') '41.576.46'
Set-CellText $ws.Range('E18') '  +0.58%  '
Set-CellText $ws.Range('D19') '6.30'
Set-CellText $ws.Range('E19') '  -1.14%  '
Set-CellText $ws.Range('D20') '0.0₃0919'
Set-CellText $ws.Range('E20') '  +0.43%  '
Set-CellText $ws.Range('D21') '70.67'
Set-CellText $ws.Range('E21') '  +2.35%  '
Set-CellText $ws.Range('D22') '11.08'
Set-CellText $ws.Range('E22') '  -4.26%  '
Set-CellText $ws.Range('D23') '235.00'
Set-CellText $ws.Range('E23') '  -0.90%  '
Set-CellText $ws.Range('E24') '  -3.07%  '
Set-CellText $ws.Range('D25') '0.999'
Set-CellText $ws.Range('E25') '  -0.22%  '
Set-CellText $ws.Range('D26') '1.88'
Set-CellText $ws.Range('E26') '  -3.35%  '
Set-CellText $ws.Range('D27') '24.37'
Set-CellText $ws.Range('E27') '  -1.58%  '
Set-CellText $ws.Range('E28') '  +0.92%  '
Set-CellText $ws.Range('D29') '9.57'
Set-CellText $ws.Range('E29') '  -1.56%  '
Set-CellText $ws.Range('D30') '36.07'
Set-CellText $ws.Range('E30') '  -1.01%  '
Set-CellText $ws.Range('D31') '153.58'
Set-CellText $ws.Range('E31') '  +0.02%  '
Set-CellText $ws.Range('D32') '5.39'
Set-CellText $ws.Range('E32') '  -4.31%  '
Set-CellText $ws.Range('D33') '2.56'
Set-CellText $ws.Range('E33') '  -2.42%  '
Set-CellText $ws.Range('E34') '  +0.35%  '
Set-CellText $ws.Range('D35') '17.92'
Set-CellText $ws.Range('E35') '  +4.34%  '
Set-CellText $ws.Range('D36') '2.52'
Set-CellText $ws.Range('E36') '  -1.37%  '
Set-CellText $ws.Range('E37') '  -1.68%  '
Set-CellText $ws.Range('D38') '1.83'
Set-CellText $ws.Range('E38') '  -3.00%  '
Set-CellText $ws.Range('E39') '  -1.39%  '
Set-CellText $ws.Range('D40') '0.100'
Set-CellText $ws.Range('E40') '  -4.44%  '
Set-CellText $ws.Range('D41') '4.10'
Set-CellText $ws.Range('E41') '  +0.67%  '
Set-CellText $ws.Range('E42') '  -0.05%  '
Set-CellText $ws.Range('D43') '19.53'
Set-CellText $ws.Range('E43') '  -9.10%  '
Set-CellText $ws.Range('D44') '1.942.20'
Set-CellText $ws.Range('E44') '  -2.20%  '
Set-CellText $ws.Range('E45') '  -0.87%  '
Set-CellText $ws.Range('D46') '2.94'
Set-CellText $ws.Range('E46') '  -3.84%  '
Set-CellText $ws.Range('D47') '8.76'
Set-CellText $ws.Range('E47') '  +0.11%  '
Set-CellText $ws.Range('D48') '2.730.10'
Set-CellText $ws.Range('E48') '  +1.07%  '
Set-CellText $ws.Range('D49') '95.62'
Set-CellText $ws.Range('E49') '  -1.95%  '
Set-CellText $ws.Range('E50') '  -2.69%  '
Set-CellText $ws.Range('D51') '66.58'
Set-CellText $ws.Range('E51') '  -4.43%  '

$excel.CutCopyMode = 0
